$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "完成对C组的测试需求说明书、测试报告、软件的评审工作"
$ws.Range("C9").Value = "基本满足要求，文档部分格式有问题"
$ws.Range("E9").Value = "王康明"
$ws.Range("F9").Value = "2h"

# Row 10
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "完成对C组的测试需求说明书、测试报告、软件的评审工作"
$ws.Range("C10").Value = "基本满足要求，软件完成情况较为完善"
$ws.Range("E10").Value = "王康明"
$ws.Range("F10").Value = "3h"

# Row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "开会确定对A、B组评审意见的接情况，根据自己负责的评审意见，对本小组的测试报告和软件，进行修改"
$ws.Range("E11").Value = "王康明"
$ws.Range("F11").Value = "2h"

# Row 12
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "开会确定对C、D组评审意见的接情况，根据自己负责的评审意见，对本小组的测试报告和软件，进行修改"
$ws.Range("E12").Value = "王康明"
$ws.Range("F12").Value = "3h"

# Row heights (match autofit result from the authored workbook)
$ws.Rows.Item(9).RowHeight = 28
$ws.Rows.Item(10).RowHeight = 28
$ws.Rows.Item(11).RowHeight = 56
$ws.Rows.Item(12).RowHeight = 56

# Restore the final selection state recorded in the workbook
$ws.Range("H12").Select() | Out-Null
